$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for every cell touched by this update (rankings, prices, % changes).
# Column D prices that are plain decimal numbers get their NumberFormat set to Text ("@")
# first so Excel stores the exact original literal (e.g. trailing zeros) instead of
# silently re-parsing the string into a floating point number.
$numericLookingCells = @(
    'D5'
    'D6'
    'D8'
    'D10'
    'D12'
    'D13'
    'D16'
    'D17'
    'D19'
    'D20'
    'D21'
    'D22'
    'D23'
    'D24'
    'D25'
    'D27'
    'D28'
    'D29'
    'D31'
    'D32'
    'D33'
    'D34'
    'D35'
    'D36'
    'D37'
    'D39'
    'D40'
    'D41'
    'D42'
    'D43'
    'D44'
    'D45'
    'D46'
    'D47'
    'D48'
    'D49'
    'D50'
    'D51'
)
foreach ($cell in $numericLookingCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$values = [ordered]@{
    'D2' = '68.223.65'
    'E2' = '  -1.90%  '
    'D3' = '2.445.20'
    'E3' = '  -1.91%  '
    'E4' = '  -0.02%  '
    'D5' = '554.03'
    'E5' = '  -2.64%  '
    'D6' = '160.02'
    'E6' = '  -2.85%  '
    'E7' = '  -0.04%  '
    'D8' = '0.497'
    'E8' = '  -2.73%  '
    'D9' = '2.444.06'
    'E9' = '  -1.88%  '
    'D10' = '0.146'
    'E10' = '  -7.76%  '
    'E11' = '  -1.42%  '
    'D12' = '0.332'
    'E12' = '  -6.16%  '
    'D13' = '4.73'
    'E13' = '  -3.66%  '
    'D14' = '2.895.00'
    'E14' = '  -1.64%  '
    'D15' = '68.197.93'
    'E15' = '  -1.68%  '
    'D16' = '0.0000165'
    'E16' = '  -5.56%  '
    'D17' = '23.04'
    'E17' = '  -5.02%  '
    'D18' = '2.458.65'
    'E18' = '  -1.01%  '
    'D19' = '10.65'
    'E19' = '  -4.57%  '
    'D20' = '338.29'
    'E20' = '  -2.12%  '
    'D21' = '6.96'
    'E21' = '  -5.43%  '
    'D22' = '3.73'
    'E22' = '  -3.67%  '
    'D23' = '0.999'
    'E23' = '  -0.05%  '
    'D24' = '1.85'
    'E24' = '  -3.75%  '
    'D25' = '65.99'
    'E25' = '  -5.20%  '
    'B26' = 'WrappedeETH'
    'C26' = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
    'D26' = '2.573.71'
    'E26' = '  -1.56%  '
    'B27' = 'NEARProtocol'
    'C27' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D27' = '3.62'
    'E27' = '  -7.28%  '
    'D28' = '0.998'
    'E28' = '  +0.50%  '
    'D29' = '7.99'
    'E29' = '  -7.76%  '
    'D30' = '0.0₃0800'
    'E30' = '  -8.52%  '
    'D31' = '7.04'
    'E31' = '  -7.91%  '
    'D32' = '0.998'
    'E32' = '  -0.10%  '
    'D33' = '425.63'
    'E33' = '  -2.57%  '
    'D34' = '1.11'
    'E34' = '  -6.21%  '
    'D35' = '1.60'
    'E35' = '  -5.87%  '
    'D36' = '155.55'
    'E36' = '  +0.42%  '
    'D37' = '18.98'
    'E37' = '  -0.39%  '
    'E38' = '  +0.00%  '
    'D39' = '0.108'
    'E39' = '  -4.23%  '
    'D40' = '17.64'
    'E40' = '  -2.81%  '
    'D41' = '0.299'
    'E41' = '  -4.76%  '
    'B42' = 'OKB'
    'C42' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'D42' = '37.29'
    'E42' = '  -1.57%  '
    'B43' = 'RenderToken'
    'C43' = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
    'D43' = '4.32'
    'E43' = '  -5.86%  '
    'B44' = 'Stacks'
    'C44' = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    'D44' = '1.44'
    'E44' = '  -8.64%  '
    'B45' = 'ImmutableX'
    'C45' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D45' = '1.08'
    'E45' = '  +1.21%  '
    'D46' = '131.07'
    'E46' = '  -5.18%  '
    'B47' = 'dogwifhat'
    'C47' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'D47' = '1.99'
    'E47' = '  -8.28%  '
    'B48' = 'Filecoin'
    'C48' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D48' = '3.29'
    'E48' = '  -4.06%  '
    'B49' = 'Cronos'
    'C49' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'D49' = '0.0710'
    'E49' = '  -1.83%  '
    'B50' = 'ARBITRUM'
    'C50' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D50' = '0.474'
    'E50' = '  -7.32%  '
    'B51' = 'Mantle'
    'C51' = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
    'D51' = '0.555'
    'E51' = '  -3.00%  '
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}